# Updates the cryptos list values (prices / 1h volume %) on Sheet1, and
# fixes the Chainlink / WrappedliquidstakedEther2.0 row ordering (rows 13-14)
# plus swaps RocketPoolETH (row 51) out for HuobiToken, matching the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: several Price values in column D are plain numeric-looking strings
# (e.g. "234.68", "0.0000101") that Excel would otherwise silently coerce
# into real numbers (and sometimes re-render in scientific notation) if
# assigned directly. Prefixing with a leading apostrophe forces Excel to
# keep them as literal text, matching the original inlineStr content.

$ws.Range("D2").Value = "43.677.77"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.351.54"
$ws.Range("E3").Value = "  +4.58%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'234.68"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").Value = "'73.60"
$ws.Range("E7").Value = "  +14.79%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.538"
$ws.Range("E9").Value = "  +23.26%  "
$ws.Range("D10").Value = "'0.0982"
$ws.Range("E10").Value = "  +3.35%  "
$ws.Range("D11").Value = "'27.63"
$ws.Range("E11").Value = "  +5.82%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.698.56"
$ws.Range("E13").Value = "  +4.35%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'16.96"
$ws.Range("E14").Value = "  +14.13%  "
$ws.Range("D15").Value = "'6.65"
$ws.Range("E15").Value = "  +10.50%  "
$ws.Range("D16").Value = "'0.887"
$ws.Range("E16").Value = "  +8.24%  "
$ws.Range("D17").Value = "2.337.17"
$ws.Range("E17").Value = "  +3.91%  "
$ws.Range("D18").Value = "43.659.68"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "'0.0000101"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").Value = "'76.11"
$ws.Range("E20").Value = "  +4.56%  "
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = "  +4.21%  "
$ws.Range("D22").Value = "'250.31"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").Value = "'3.81"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +3.01%  "
$ws.Range("D26").Value = "'10.25"
$ws.Range("E26").Value = "  +6.16%  "
$ws.Range("D27").Value = "'2.25"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("E28").Value = "  +3.98%  "
$ws.Range("D29").Value = "'172.41"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("E30").Value = "  +7.98%  "
$ws.Range("D31").Value = "'0.132"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("E32").Value = "  +4.56%  "
$ws.Range("D33").Value = "'5.11"
$ws.Range("E33").Value = "  +4.42%  "
$ws.Range("E34").Value = "  +3.28%  "
$ws.Range("D35").Value = "'5.09"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").Value = "'3.76"
$ws.Range("E36").Value = "  +3.90%  "
$ws.Range("D37").Value = "'2.42"
$ws.Range("E37").Value = "  +7.37%  "
$ws.Range("D38").Value = "'6.37"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +6.20%  "
$ws.Range("D40").Value = "'19.48"
$ws.Range("E40").Value = "  +14.38%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +0.98%  "
$ws.Range("E43").Value = "  +8.86%  "
$ws.Range("D44").Value = "'1.22"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").Value = "'98.69"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("D46").Value = "'0.0961"
$ws.Range("E46").Value = "  +2.61%  "
$ws.Range("D47").Value = "'4.43"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'0.180"
$ws.Range("E48").Value = "  +13.31%  "
$ws.Range("D49").Value = "1.437.96"
$ws.Range("E49").Value = "  +1.12%  "
$ws.Range("D50").Value = "'2.28"
$ws.Range("E50").Value = "  +1.73%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.77"
$ws.Range("E51").Value = "  +2.01%  "
